# VSSC Feedback workbook update
# - Follow Up sheet becomes the active tab (was Interim... wait see below)
# - Several "GG0130A" feedback rows get a checkmark / placeholder answer column filled in
# - Two new remark strings are introduced

$wb = $excel.ActiveWorkbook

$wsInterim    = $wb.Worksheets.Item("Interim")
$wsDischarge  = $wb.Worksheets.Item("Discharge")
$wsFollowUp   = $wb.Worksheets.Item("Follow Up")

# ---------------------------------------------------------------------------
# Discharge sheet ("Discharge" tab) cell updates
# ---------------------------------------------------------------------------

# Row 8: new answer cell C8 = "????" (inherits the row's style automatically)
$wsDischarge.Range("C8").Value = "????"

# Row 10: C10 changes from the "in conflict with..." note to a checkmark,
# and picks up the checkmark column's style (copy format from C2, a cell
# that already carries that style, then set the value)
$wsDischarge.Range("C2").Copy()
$wsDischarge.Range("C10").PasteSpecial(-4122)
$wsDischarge.Range("C10").Value = [char]0x221A

# Row 19: C19 changes from "????" to a checkmark and adopts the checkmark style
$wsDischarge.Range("C2").Copy()
$wsDischarge.Range("C19").PasteSpecial(-4122)
$wsDischarge.Range("C19").Value = [char]0x221A

# Row 23: C23 changes from a checkmark to the new remark
# (must be written AFTER the Follow Up sheet's "Assessment Completed" cell
# below so the shared-string table gets the two new strings in the right order)

# Row 25: C25 changes from "????" to a checkmark (style stays the same)
$wsDischarge.Range("C25").Value = [char]0x221A

# ---------------------------------------------------------------------------
# Follow Up sheet cell updates
# ---------------------------------------------------------------------------

# Row 19: C19 changes from "????" to a checkmark, adopting the checkmark style
$wsFollowUp.Range("C17").Copy()
$wsFollowUp.Range("C19").PasteSpecial(-4122)
$wsFollowUp.Range("C19").Value = [char]0x221A

# Row 23: new answer cell C23 = checkmark, adopting the checkmark style
$wsFollowUp.Range("C17").Copy()
$wsFollowUp.Range("C23").PasteSpecial(-4122)
$wsFollowUp.Range("C23").Value = [char]0x221A

# Row 35: new answer cell C35 = "Assessment Completed" (new shared string,
# must be entered before the Discharge sheet's "Discharge Performance and
# Discharge Goal?" string so the shared-string table order matches)
$wsFollowUp.Range("C35").Value = "Assessment Completed"

# Now that "Assessment Completed" exists, add the Discharge sheet's new remark
$wsDischarge.Range("C23").Value = "Discharge Performance and Discharge Goal?"

# ---------------------------------------------------------------------------
# Selections / active sheet
# ---------------------------------------------------------------------------

# Update the remembered selections on Discharge and Follow Up before we
# switch the active tab away from them.
$wsDischarge.Range("D14").Select()
$wsFollowUp.Range("C35").Select()

# The "Interim" tab becomes the selected/active sheet (was "Follow Up").
$wsInterim.Activate()
